# cryptos.xlsx refresh — GitHub Actions bot run (Mon Mar 27 00:55:58 UTC 2023)
#
# Re-pulls the coinranking.com snapshot into Sheet1: every coin's Price (col D)
# and 1h Volume change (col E) is refreshed to the latest quote, and the two
# rank swaps that happened this run (Stellar now just edges out Filecoin, and
# VeChain now just edges out Hedera) are reflected by writing the new coin's
# Name/Link/Price/Volume into the row that kept its rank position.
#
# All data cells on this sheet are plain text (prices like "28.042.08" or
# "1.001" are display strings, not numbers — note the thousands-dot grouping),
# so every write below forces the Text number format before assigning the
# value and restores the cell's original Style afterwards. That keeps Excel
# from "helpfully" reinterpreting a string like "1.001" as the number 1.001,
# while leaving formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target A1-style cell -> new text value, in row order (Coin/Link/Price/Volume
# columns B-E; only the cells that actually changed this run are listed).
$updates = @(
    @{ Cell = 'D2'; Value = '28.032.69' },
    @{ Cell = 'E2'; Value = '  +1.60%  ' },
    @{ Cell = 'D3'; Value = '1.780.88' },
    @{ Cell = 'E3'; Value = '  +1.51%  ' },
    @{ Cell = 'D4'; Value = '1.001' },
    @{ Cell = 'E4'; Value = '  +0.10%  ' },
    @{ Cell = 'D5'; Value = '329.13' },
    @{ Cell = 'E5'; Value = '  +1.85%  ' },
    @{ Cell = 'D6'; Value = '1.000' },
    @{ Cell = 'E6'; Value = '  +0.11%  ' },
    @{ Cell = 'D7'; Value = '0.4496' },
    @{ Cell = 'E7'; Value = '  +0.12%  ' },
    @{ Cell = 'D8'; Value = '0.3551' },
    @{ Cell = 'E8'; Value = '  +0.70%  ' },
    @{ Cell = 'D9'; Value = '0.07451' },
    @{ Cell = 'E9'; Value = '  +0.49%  ' },
    @{ Cell = 'D10'; Value = '42.15' },
    @{ Cell = 'E10'; Value = '  +0.99%  ' },
    @{ Cell = 'D11'; Value = '1.106' },
    @{ Cell = 'E11'; Value = '  +1.84%  ' },
    @{ Cell = 'D12'; Value = '1.001' },
    @{ Cell = 'E12'; Value = '  +0.20%  ' },
    @{ Cell = 'D13'; Value = '20.93' },
    @{ Cell = 'E13'; Value = '  +1.36%  ' },
    @{ Cell = 'D14'; Value = '6.046' },
    @{ Cell = 'E14'; Value = '  +1.50%  ' },
    @{ Cell = 'D15'; Value = '7.251' },
    @{ Cell = 'E15'; Value = '  +1.62%  ' },
    @{ Cell = 'D16'; Value = '1.775.21' },
    @{ Cell = 'E16'; Value = '  +1.52%  ' },
    @{ Cell = 'D17'; Value = '93.50' },
    @{ Cell = 'E17'; Value = '  +1.90%  ' },
    @{ Cell = 'D18'; Value = '0.00001063' },
    @{ Cell = 'E18'; Value = '  +0.56%  ' },
    @{ Cell = 'D19'; Value = '0.06416' },
    @{ Cell = 'E19'; Value = '  +0.15%  ' },
    @{ Cell = 'D20'; Value = '1.001' },
    @{ Cell = 'E20'; Value = '  +0.18%  ' },
    @{ Cell = 'D21'; Value = '17.13' },
    @{ Cell = 'E21'; Value = '  +0.50%  ' },
    @{ Cell = 'D22'; Value = '5.791' },
    @{ Cell = 'E22'; Value = '  +0.84%  ' },
    @{ Cell = 'D23'; Value = '28.080.29' },
    @{ Cell = 'E23'; Value = '  +1.65%  ' },
    @{ Cell = 'D24'; Value = '11.33' },
    @{ Cell = 'E24'; Value = '  +1.33%  ' },
    @{ Cell = 'D25'; Value = '2.116' },
    @{ Cell = 'E25'; Value = '  -0.01%  ' },
    @{ Cell = 'D26'; Value = '161.01' },
    @{ Cell = 'E26'; Value = '  -0.46%  ' },
    @{ Cell = 'D27'; Value = '20.30' },
    @{ Cell = 'E27'; Value = '  +0.28%  ' },
    @{ Cell = 'D28'; Value = '1.988.28' },
    @{ Cell = 'E28'; Value = '  +2.02%  ' },
    @{ Cell = 'D29'; Value = '2.165' },
    @{ Cell = 'E29'; Value = '  +5.42%  ' },
    @{ Cell = 'D30'; Value = '124.59' },
    @{ Cell = 'E30'; Value = '  -0.81%  ' },
    @{ Cell = 'D31'; Value = '1.101' },
    @{ Cell = 'E31'; Value = '  +4.38%  ' },
    @{ Cell = 'B32'; Value = 'Stellar' },
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' },
    @{ Cell = 'D32'; Value = '0.09211' },
    @{ Cell = 'E32'; Value = '  +0.81%  ' },
    @{ Cell = 'B33'; Value = 'Filecoin' },
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' },
    @{ Cell = 'D33'; Value = '5.682' },
    @{ Cell = 'E33'; Value = '  +4.05%  ' },
    @{ Cell = 'D34'; Value = '3.679' },
    @{ Cell = 'E34'; Value = '  +0.29%  ' },
    @{ Cell = 'D35'; Value = '11.93' },
    @{ Cell = 'E35'; Value = '  +2.32%  ' },
    @{ Cell = 'B36'; Value = 'VeChain' },
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Cell = 'D36'; Value = '0.02284' },
    @{ Cell = 'E36'; Value = '  +0.18%  ' },
    @{ Cell = 'B37'; Value = 'Hedera' },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' },
    @{ Cell = 'D37'; Value = '0.06152' },
    @{ Cell = 'E37'; Value = '  +1.87%  ' },
    @{ Cell = 'D38'; Value = '0.2105' },
    @{ Cell = 'E38'; Value = '  +1.80%  ' },
    @{ Cell = 'D39'; Value = '4.987' },
    @{ Cell = 'E39'; Value = '  +1.20%  ' },
    @{ Cell = 'D40'; Value = '0.6321' },
    @{ Cell = 'E40'; Value = '  +0.80%  ' },
    @{ Cell = 'D41'; Value = '1.182' },
    @{ Cell = 'E41'; Value = '  +0.10%  ' },
    @{ Cell = 'D42'; Value = '1.383' },
    @{ Cell = 'E42'; Value = '  -0.20%  ' },
    @{ Cell = 'D43'; Value = '7.890' },
    @{ Cell = 'E43'; Value = '  +2.01%  ' },
    @{ Cell = 'D44'; Value = '13.25' },
    @{ Cell = 'E44'; Value = '  +0.78%  ' },
    @{ Cell = 'D45'; Value = '3.748' },
    @{ Cell = 'E45'; Value = '  +1.13%  ' },
    @{ Cell = 'D46'; Value = '0.5904' },
    @{ Cell = 'E46'; Value = '  +0.79%  ' },
    @{ Cell = 'D47'; Value = '123.04' },
    @{ Cell = 'E47'; Value = '  +0.35%  ' },
    @{ Cell = 'D48'; Value = '1.957' },
    @{ Cell = 'E48'; Value = '  +1.00%  ' },
    @{ Cell = 'D49'; Value = '1.141' },
    @{ Cell = 'E49'; Value = '  +1.88%  ' },
    @{ Cell = 'D50'; Value = '0.06893' },
    @{ Cell = 'E50'; Value = '  -0.59%  ' },
    @{ Cell = 'B51'; Value = 'Tezos' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/fsIbGOEJWbzxG+tezos-xtz' },
    @{ Cell = 'D51'; Value = '1.147' },
    @{ Cell = 'E51'; Value = '  +3.04%  ' }
)

foreach ($update in $updates) {
    $cellRef = $update.Cell
    $value = $update.Value
    $range = $ws.Range($cellRef)

    # Column D holds prices; some are plain decimals ("1.001", "93.50") that
    # Excel's General format would silently convert to numbers. Anything else
    # (coin names, links, "28.042.08"-style multi-dot prices, the "  +x.xx%  "
    # volume strings) is never mistaken for a number, so it's safe as-is.
    $looksNumeric = $cellRef -match '^D\d+$' -and $value -match '^[+-]?\d+(\.\d+)?$'

    if ($looksNumeric) {
        $originalStyle = $range.Style
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = $originalStyle
    } else {
        $range.Value = $value
    }
}
